$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a number need NumberFormat forced to "Text" so
# Excel keeps storing a string (matching the source inlineStr cells) instead of
# silently parsing it into a numeric value. The Style reset afterwards drops the
# temporary text-format style again so the cell stays on the default style, exactly
# like the untouched cells around it.

$ws.Range("D2").Value = "42.864.77"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.325.07"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "2.685.93"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "2.391.77"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").Value = "42.818.47"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.04%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "0.0₃0891"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  +7.53%  "
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("E30").Value = "  -2.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "145.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("E40").Value = "  +0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +22.76%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "1.927.41"
$ws.Range("E43").Value = "  -3.45%  "
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "2.554.34"
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
